$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet
$ws.Name = "Redemptions"

# 2. Update reporting period text (shared string used by C4 / becomes B4)
$ws.Range("C4").Value = "01-Jul-2023 To 24-Jun-2024"

# 5. Shift header/summary block from columns B/C to A/B (Copy preserves style+format)
$ws.Range("B1").Copy($ws.Range("A1"))
$ws.Range("B1").Clear()

$ws.Range("B2").Copy($ws.Range("A2"))
$ws.Range("C2").Copy($ws.Range("B2"))
$ws.Range("C2").Clear()

$ws.Range("B3").Copy($ws.Range("A3"))
$ws.Range("C3").Copy($ws.Range("B3"))
$ws.Range("C3").Clear()

$ws.Range("B4").Copy($ws.Range("A4"))
$ws.Range("C4").Copy($ws.Range("B4"))
$ws.Range("C4").Clear()

$ws.Range("B5").Copy($ws.Range("A5"))
$ws.Range("C5").Copy($ws.Range("B5"))
$ws.Range("C5").Clear()

$ws.Range("B6").Copy($ws.Range("A6"))
$ws.Range("C6").Copy($ws.Range("B6"))
$ws.Range("C6").Clear()

# 3. Left-align the value now in B5:B6 (formerly C5:C6, style index 6)
$ws.Range("B5:B6").HorizontalAlignment = -4131

# 4. Adjust column widths (closest achievable values given engine's pixel-width quantization)
$ws.Columns.Item(1).ColumnWidth = 30
$ws.Columns.Item(2).ColumnWidth = 26.5
$ws.Columns.Item(3).ColumnWidth = 8.33
